$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for each coin row
# to reflect the refreshed crypto data from the GitHub Actions run.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.394.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.61%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.109.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.10%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.50%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.40%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5237'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.57%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4451'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.42%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.97%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09383'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.68%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.176'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.12%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.75%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.696'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.68%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.111.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.90%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.933'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.98%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.12%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001164'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.79%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.005'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.39%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.42%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06725'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.48%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.310'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.28%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.005'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.426.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.55%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.32%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.319'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.23%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.385.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.08%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.93%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.542'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.08%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.85%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.150'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.18%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.763'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.41%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1058'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.90%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.826'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.39%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.269'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.52%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.942'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.51%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.77%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02640'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.89%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06840'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.06%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7080'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.89%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.31%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.338'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.94%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2237'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.33%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6864'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.91%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.60'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.72%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.376'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.88%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.005'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.28%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.401'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +20.00%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.652'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.28%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.224'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.73%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.202'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.50%  '

Write-Output "Cryptos list updated"